$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transition_Matrix")
$rng = $ws.Range("A1:D49")
$key1 = $ws.Range("A2")
$rng.Sort($key1, 1, $null, $null, 1, $null, 1, 1)
for ($r = 1; $r -le 49; $r++) {
  Write-Host "$r : $($ws.Cells.Item($r,1).Value2) | $($ws.Cells.Item($r,2).Value2) | $($ws.Cells.Item($r,3).Value2) | $($ws.Cells.Item($r,4).Value2)"
}
